$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "MCT-3A-Processos de Usinagem 1"
$ws.Range("D2").Value = "-"

# Row 3
$ws.Range("B3").Value = "MCT-3A-Processos de Usinagem 1"
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("B4").Value = "MCT-3A-Processos de Usinagem 1"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "MEC-2B-Tornearia"

# Row 6
$ws.Range("B6").Value = "MCT-3A-Processos de Usinagem 1"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "MEC-2B-Tornearia"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "MEC-2B-Tornearia"

# Row 8
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "MEC-2B-Tornearia"

# Row 12
$ws.Range("B12").Value = "MEC-2A-Tornearia"
$ws.Range("F12").Value = "-"

# Row 14
$ws.Range("B14").Value = "MEC-2A-Tornearia"
$ws.Range("F14").Value = "-"

# Row 15
$ws.Range("B15").Value = "MEC-2A-Tornearia"
$ws.Range("F15").Value = "-"

# Row 16
$ws.Range("B16").Value = "MEC-2A-Tornearia"
$ws.Range("F16").Value = "-"

$wb.Save()
